$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.488.58'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '1.813.32'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = '''225.91'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('E6').Value = '  +2.82%  '
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').Value = '''38.39'
$ws.Range('E8').Value = '  +7.11%  '
$ws.Range('E9').Value = '  -4.04%  '
$ws.Range('E10').Value = '  -2.75%  '
$ws.Range('E11').Value = '  +0.86%  '
$ws.Range('D12').Value = '2.073.69'
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('D13').Value = '''11.20'
$ws.Range('E13').Value = '  -2.05%  '
$ws.Range('D14').Value = '1.822.10'
$ws.Range('E15').Value = '  -1.87%  '
$ws.Range('D16').Value = '34.465.69'
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('E17').Value = '  -2.14%  '
$ws.Range('D18').Value = '''68.33'
$ws.Range('E18').Value = '  -1.33%  '
$ws.Range('D19').Value = '''242.63'
$ws.Range('E19').Value = '  -1.33%  '
$ws.Range('E20').Value = '  -2.89%  '
$ws.Range('E21').Value = '  -2.20%  '
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('D23').Value = '''4.13'
$ws.Range('E23').Value = '  -1.68%  '
$ws.Range('E24').Value = '  +3.72%  '
$ws.Range('D25').Value = '''170.23'
$ws.Range('E25').Value = '  -0.38%  '
$ws.Range('D26').Value = '''7.83'
$ws.Range('E26').Value = '  -1.31%  '
$ws.Range('D27').Value = '''17.59'
$ws.Range('E27').Value = '  +3.69%  '
$ws.Range('E28').Value = '  +1.59%  '
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('E30').Value = '  -1.72%  '
$ws.Range('E31').Value = '  -1.41%  '
$ws.Range('E32').Value = '  -2.71%  '
$ws.Range('E33').Value = '  -5.80%  '
$ws.Range('E34').Value = '  -0.51%  '
$ws.Range('D35').Value = '1.364.17'
$ws.Range('E35').Value = '  -2.44%  '
$ws.Range('D36').Value = '''0.646'
$ws.Range('E36').Value = '  -4.07%  '
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('E38').Value = '  -5.56%  '
$ws.Range('E39').Value = '  -1.72%  '
$ws.Range('E40').Value = '  -0.86%  '
$ws.Range('E41').Value = '  +1.23%  '
$ws.Range('E42').Value = '  -1.42%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '''2.80'
$ws.Range('E43').Value = '  -0.84%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '''81.65'
$ws.Range('E44').Value = '  -1.25%  '
$ws.Range('D45').Value = '''13.81'
$ws.Range('E45').Value = '  +2.42%  '
$ws.Range('E46').Value = '  +1.52%  '
$ws.Range('D47').Value = '1.974.88'
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('E48').Value = '  -4.49%  '
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('D50').Value = '''102.05'
$ws.Range('E50').Value = '  -3.04%  '
$ws.Range('E51').Value = '  -5.16%  '
